$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 5678
$wsExhibit.Range("G3").Value = 75
$wsExhibit.Range("F6").Value = 951
$wsExhibit.Range("F8").Value = 2560
$wsExhibit.Range("F10").Value = 171
$wsExhibit.Range("F12").Value = 87
$wsExhibit.Range("F14").Value = 2396
$wsExhibit.Range("F15").Value = 413

# --- Sheet "演出" (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("G2").Value = 98

# --- Sheet "全部类型" (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 5678
$wsAll.Range("G3").Value = 75
$wsAll.Range("G4").Value = 98
$wsAll.Range("F8").Value = 951
$wsAll.Range("F10").Value = 2560
$wsAll.Range("F12").Value = 171
$wsAll.Range("F15").Value = 87
$wsAll.Range("F17").Value = 2396
$wsAll.Range("F18").Value = 413
